# Update Approved/Rejected + ReasonToReject for TestScenario_8 rows (8 & 9)
# and move the active selection to H8, matching the authored change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

$ws.Range("I8").Value = "Rejected"
$ws.Range("J8").Value = "Nil"

$ws.Range("I9").Value = "Rejected"
$ws.Range("J9").Value = "Nil"

$ws.Range("H8").Select() | Out-Null
